$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook level: rename "Sheet1" -> "roles" and make "photos" the active
#    tab (the diff bumps workbookView activeTab from 3 to 4, zero-based, and
#    sheet5 / photos becomes tabSelected="1").
# ---------------------------------------------------------------------------
$wsRoles  = $wb.Worksheets.Item("Sheet1")
$wsRoles.Name = "roles"

$wsRights  = $wb.Worksheets.Item("rights")
$wsGrants  = $wb.Worksheets.Item("grants")
$wsPhotos  = $wb.Worksheets.Item("photos")

# ---------------------------------------------------------------------------
# 2. "roles" sheet (sheet4) — add the SQL-generating helper column (F) for
#    the role_id = 4 (photographer) block, drop the stray "members / READ"
#    row that used to sit at row 20, add the "member: role_id = 5" section
#    header, and add the role_id = 5 (member) grants block with its own
#    helper column.
# ---------------------------------------------------------------------------

# role_id = 4 block (rows 4-15) helper formulas
$wsRoles.Range("F4").Formula = '="INSERT INTO `grants` (`right_id`,`role_id`,`created_at`,`updated_at`) VALUES ("&B4&",4,''2012-05-22'',''2012-05-22'');"'
$wsRoles.Range("F5:F15").Formula = '="INSERT INTO `grants` (`right_id`,`role_id`,`created_at`,`updated_at`) VALUES ("&B5&",4,''2012-05-22'',''2012-05-22'');"'

# Remove the old "members / READ" row that lived at row 20.
$wsRoles.Range("B20:D20").ClearContents()

# New section header for the member role.
$wsRoles.Range("A19").Value = "member: role_id = 5"

# role_id = 5 block (rows 21-24 already existed; 25-29 are new)
$wsRoles.Range("B25").Value = 92
$wsRoles.Range("C25").Value = "users"
$wsRoles.Range("D25").Value = "UPDATE"

$wsRoles.Range("B26").Value = 91
$wsRoles.Range("C26").Value = "users"
$wsRoles.Range("D26").Value = "READ"

$wsRoles.Range("B27").Value = 25
$wsRoles.Range("C27").Value = "members"
$wsRoles.Range("D27").Value = "READ"

$wsRoles.Range("B28").Value = 26
$wsRoles.Range("C28").Value = "members"
$wsRoles.Range("D28").Value = "UPDATE"

$wsRoles.Range("B29").Value = 28
$wsRoles.Range("C29").Value = "pages"
$wsRoles.Range("D29").Value = "READ"

# role_id = 5 helper formulas (rows 21-29)
$wsRoles.Range("F21").Formula = '="INSERT INTO `grants` (`right_id`,`role_id`,`created_at`,`updated_at`) VALUES ("&B21&",5,''2012-05-22'',''2012-05-22'');"'
$wsRoles.Range("F22:F29").Formula = '="INSERT INTO `grants` (`right_id`,`role_id`,`created_at`,`updated_at`) VALUES ("&B22&",5,''2012-05-22'',''2012-05-22'');"'

$wsRoles.Range("D32").Select()

# ---------------------------------------------------------------------------
# 3. "photos" sheet (sheet5) — append the 41 "Annual Repairs, Etc." (2012
#    airshow maintenance gallery) photo rows (30-70): filename in F, order
#    in G, generated INSERT statement in K.
# ---------------------------------------------------------------------------
$filenames = @(
  "Annual_Repairs_Etc_2012_01.jpg","Annual_Repairs_Etc_2012_02.jpg","Annual_Repairs_Etc_2012_03.jpg",
  "Annual_Repairs_Etc_2012_04.jpg","Annual_Repairs_Etc_2012_05.jpg","Annual_Repairs_Etc_2012_06.jpg",
  "Annual_Repairs_Etc_2012_07.jpg","Annual_Repairs_Etc_2012_08.jpg","Annual_Repairs_Etc_2012_09.jpg",
  "Annual_Repairs_Etc_2012_10.jpg","Annual_Repairs_Etc_2012_11.jpg","Annual_Repairs_Etc_2012_12.jpg",
  "Annual_Repairs_Etc_2012_13.jpg","Annual_Repairs_Etc_2012_14.jpg","Annual_Repairs_Etc_2012_15.jpg",
  "Annual_Repairs_Etc_2012_16.jpg","Annual_Repairs_Etc_2012_17.jpg","Annual_Repairs_Etc_2012_18.jpg",
  "Annual_Repairs_Etc_2012_19.jpg","Annual_Repairs_Etc_2012_20.jpg","Annual_Repairs_Etc_2012_21.jpg",
  "Annual_Repairs_Etc_2012_22.jpg","Annual_Repairs_Etc_2012_23.jpg","Annual_Repairs_Etc_2012_24.jpg",
  "Annual_Repairs_Etc_2012_25.jpg","Annual_Repairs_Etc_2012_26.jpg","Annual_Repairs_Etc_2012_27.jpg",
  "Annual_Repairs_Etc_2012_28.jpg","Annual_Repairs_Etc_2012_29.jpg","Annual_Repairs_Etc_2012_30.jpg",
  "Annual_Repairs_Etc_2012_31.jpg","Annual_Repairs_Etc_2012_32.jpg","Annual_Repairs_Etc_2012_33.jpg",
  "Annual_Repairs_Etc_2012_34.jpg","Annual_Repairs_Etc_2012_35.jpg","Annual_Repairs_Etc_2012_36.jpg",
  "Annual_Repairs_Etc_2012_37.jpg","Annual_Repairs_Etc_2012_38.jpg","Annual_Repairs_Etc_2012_39.jpg",
  "Annual_Repairs_Etc_2012_40.jpg","Annual_Repairs_Etc_2012_41.jpg"
)

for ($i = 0; $i -lt $filenames.Length; $i++) {
    $row = 30 + $i
    $order = $i + 1
    $wsPhotos.Cells.Item($row, 6).Value = $filenames[$i]
    $wsPhotos.Cells.Item($row, 7).Value = $order
}

$wsPhotos.Range("K30").Formula = '="INSERT INTO `photos` (`gallery_id`,`photographer_id`,`filename`,`order`,`created_at`,`updated_at`,`path`) VALUES (15,9,''"&F30&"'',''"&G30&"'',''2012-02-28'',''2012-02-28'',''2012_maintenance'');"'
$wsPhotos.Range("K31:K70").Formula = '="INSERT INTO `photos` (`gallery_id`,`photographer_id`,`filename`,`order`,`created_at`,`updated_at`,`path`) VALUES (15,9,''"&F31&"'',''"&G31&"'',''2012-02-28'',''2012-02-28'',''2012_maintenance'');"'

$wsPhotos.Activate()
$wsPhotos.Range("H30").Select()

# ---------------------------------------------------------------------------
# 4. View-state tweaks on "rights" and "grants" captured by the diff (pure
#    scroll/selection changes — no data changed on those two sheets).
# ---------------------------------------------------------------------------
$wsRights.Range("B29").Select()
$wsRights.Range("A48:C48").Select()

$wsGrants.Range("B38").Select()
$wsGrants.Range("C67:C93").Select()
